$d = $word.ActiveDocument

# --- 1. Locate the sentence to edit and split it where the new bookmark
#        needs to land: right after "H60L03" (soon to be "H60L04"), before
#        " with a MVC project called ".
$rng = $d.Content
$found = $rng.Find.Execute("Create a solution called H60L03", $true, $false, `
                            $false, $false, $false, $true, 1, $false, "", 0)

$splitPoint = $rng.Duplicate
$splitPoint.Collapse(0)  # wdCollapseEnd -> position right after "H60L03"

# Adding a bookmark named "_GoBack" automatically relocates the document's
# single hidden "_GoBack" bookmark here (Word only ever keeps one), which
# also splits the run in two at the insertion point.
$d.Bookmarks.Add("_GoBack", $splitPoint)

# --- 2. Rename the lab number in the (now isolated) first run.
$rng2 = $d.Content
$rng2.Find.Execute("H60L03", $true, $false, $false, $false, $false, $true, `
                    1, $false, "", 0)
$rng2.Text = "H60L04"
